# Houston roster: the underlying player records got reshuffled onto different
# row numbers (rows 5/6 swapped; rows 8/9/10 cyclically rotated; rows 14/15/16
# cyclically rotated), while the "No." column (A) stays anchored to the row
# and each row's hyperlink relationship (set up for column K) also stays
# anchored to the row. Only the textual content of columns B..K is updated to
# reflect which player's record now lives on that row. (Because the hyperlink
# target itself is not touched, only the "K" display text, the end result
# intentionally keeps each row's original bbref hyperlink target while
# displaying the new player's url as text -- matching the source edit.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $name, $jersey, $pos, $ht, $wt, $bday, $nat, $exp, $college, $urlText) {
    $ws.Cells.Item($r, 3).Value = $name          # C: Player
    $ws.Cells.Item($r, 2).Value = $jersey        # B: Jersey No.
    $ws.Cells.Item($r, 4).Value = $pos           # D: Pos
    $ws.Cells.Item($r, 5).Value = $ht            # E: Ht
    $ws.Cells.Item($r, 6).Value = $wt            # F: Wt
    $ws.Cells.Item($r, 7).Value = $bday          # G: Birth Date
    $ws.Cells.Item($r, 8).Value = $nat           # H: nationality
    $ws.Cells.Item($r, 9).Value = $exp           # I: Exp
    if ($college -eq $null) {
        $ws.Cells.Item($r, 10).ClearContents() | Out-Null
    } else {
        $ws.Cells.Item($r, 10).Value = $college   # J: College
    }
    $ws.Cells.Item($r, 11).Value = $urlText      # K: bbref url (display text only)
}

# Row 5 <- old row 6 (Alperen Sengun)
Set-Row 5 "Alperen Şengün" 28 "C" "6-9" 235 "July 25, 2002" "tr" "1" $null "https://www.basketball-reference.com/players/s/sengual01.html"

# Row 6 <- old row 5 (Jalen Green)
Set-Row 6 "Jalen Green" 4 "SG" "6-4" 178 "February 9, 2002" "us" "1" $null "https://www.basketball-reference.com/players/g/greenja05.html"

# Row 8 <- old row 10 (Josh Christopher)
Set-Row 8 "Josh Christopher" 9 "SG" "6-5" 215 "December 8, 2001" "us" "1" "Arizona State" "https://www.basketball-reference.com/players/c/chrisjo01.html"

# Row 9 <- old row 8 (Daishen Nix)
Set-Row 9 "Daishen Nix" 15 "PG" "6-5" 224 "February 13, 2002" "us" "1" $null "https://www.basketball-reference.com/players/n/nixda01.html"

# Row 10 <- old row 9 (Kevin Porter Jr.)
Set-Row 10 "Kevin Porter Jr." 3 "PG" "6-4" 203 "May 4, 2000" "us" "3" "USC" "https://www.basketball-reference.com/players/p/porteke02.html"

# Row 14 <- old row 16 (Frank Kaminsky)
Set-Row 14 "Frank Kaminsky" 33 "C" "7-0" 240 "April 4, 1993" "us" "7" "Wisconsin" "https://www.basketball-reference.com/players/k/kaminfr01.html"

# Row 15 <- old row 14 (Trevor Hudgins (TW))
Set-Row 15 "Trevor Hudgins (TW)" 12 "PG" "6-0" 180 "March 23, 1999" "us" "R" "Northwest Missouri State University" "https://www.basketball-reference.com/players/h/hudgitr01.html"

# Row 16 <- old row 15 (Darius Days (TW))
Set-Row 16 "Darius Days (TW)" 5 "PF" "6-7" 245 "October 20, 1999" "us" "R" "LSU" "https://www.basketball-reference.com/players/d/daysda01.html"
